$d = $word.ActiveDocument

# 1. Remove the standalone "Meta description: ..." paragraph that
#    currently follows the title (Heading1) paragraph.
$d.Paragraphs(2).Range.Delete()

# 2. Insert a new bold paragraph ("Play 7 Diamond and Win Big | Free Slot
#    Game") right before the last paragraph in the document.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertPoint.InsertBefore("Play 7 Diamond and Win Big | Free Slot Game" + [char]13)

$newPara = $d.Paragraphs($n)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newRange.Font.Italic = $false
$newRange.Font.Bold = $true

# 3. Replace the old "feature image" prompt text of what is now the final
#    paragraph with the review meta description, keeping its italic run
#    formatting intact.
$d.Content.Find.Execute(
  "Please create a feature image for the game ""7 Diamond"". The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be surrounded by diamonds and other symbols from the game, such as bells, fruit, and the number 7. The image should convey excitement and the possibility of winning big.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Read our review of 7 Diamond, a traditional slot game with easy gameplay, available to play free on desktop, mobile and tablet.",
  2)
